# Update "想去人数" (F column) counts that changed between crawls.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 673
$wsExpo.Range("F6").Value = 1434
$wsExpo.Range("F8").Value = 1698
$wsExpo.Range("F16").Value = 69
$wsExpo.Range("F20").Value = 59
$wsExpo.Range("F25").Value = 91

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 68
$wsShow.Range("F3").Value = 1

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 673
$wsAll.Range("F6").Value = 1434
$wsAll.Range("F8").Value = 1698
$wsAll.Range("F16").Value = 69
$wsAll.Range("F20").Value = 59
$wsAll.Range("F23").Value = 68
$wsAll.Range("F25").Value = 1
$wsAll.Range("F27").Value = 91
